# Settings Add departments screen script started and completed the add departments script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell
$ws.Range("B1").Value = "Departments"

# New data cell with wrapped text
$ws.Range("B2").Value = "AutomatonTest"
$ws.Range("B2").WrapText = $true

# Row 2 needs to grow to fit the wrapped text (matches ht="28.8" in the target sheet)
$ws.Rows.Item(2).RowHeight = 28.8

# Update the selection to mirror the authored workbook state
$ws.Range("G7").Select()
